# Insert a new data row at row 142 (weekly Fruta/Hortaliza update).
# Excel shifts the existing rows 142:152 down to 143:153 and updates the
# sheet dimension automatically, matching the canonical diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("142:142").Insert()

$row = 142

$ws.Cells.Item($row, 1).Value  = 8
$ws.Cells.Item($row, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value  = "Coquimbo"
$ws.Cells.Item($row, 4).Value  = 44946
$ws.Cells.Item($row, 5).Value  = 4
$ws.Cells.Item($row, 6).Value  = 100112052
$ws.Cells.Item($row, 7).Value  = "Albahaca"
$ws.Cells.Item($row, 8).Value  = "Sin especificar"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 700
$ws.Cells.Item($row, 11).Value = 4000
$ws.Cells.Item($row, 12).Value = 5000
$ws.Cells.Item($row, 13).Value = 4500
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 750
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
